# Update "want to go" (想去人数) counts in the 上海-漫展信息 workbook,
# mirroring the values regenerated by the gh-pages build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 289
$ws1.Range("F15").Value = 635
$ws1.Range("F16").Value = 1418
$ws1.Range("F17").Value = 1418
$ws1.Range("F23").Value = 479
$ws1.Range("F24").Value = 25747
$ws1.Range("F25").Value = 25747
$ws1.Range("F28").Value = 16326
$ws1.Range("F29").Value = 16326
$ws1.Range("F30").Value = 376
$ws1.Range("F34").Value = 188
$ws1.Range("F36").Value = 442
$ws1.Range("F39").Value = 651
$ws1.Range("F40").Value = 386

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 58
$ws2.Range("F17").Value = 384
$ws2.Range("F34").Value = 852
$ws2.Range("F43").Value = 788

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 289
$ws4.Range("F22").Value = 635
$ws4.Range("F23").Value = 1418
$ws4.Range("F24").Value = 1418
$ws4.Range("F27").Value = 58
$ws4.Range("F29").Value = 384
$ws4.Range("F31").Value = 479
$ws4.Range("F33").Value = 25747
$ws4.Range("F36").Value = 16326
$ws4.Range("F37").Value = 376
$ws4.Range("F40").Value = 188
$ws4.Range("F43").Value = 442
$ws4.Range("F46").Value = 651
